$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Apple iPhone X (Silver, 64 GB)"
$ws.Range("B1").Value = "[[ChromeDriver: chrome on XP (76f6c00741495541763e2fac4bf6cee6)] -> xpath: //div[@class='_1vC4OE _2rQ-NK']]"
